$d = $word.ActiveDocument

# The trailing "_GoBack" bookmark currently sits right after the
# "Sous-total " run (end of the paragraph list that starts with
# "Actions disponibles" ... "Cette page permet ..."). The edit moves
# that bookmark down to a brand-new final paragraph ("Tets pour la ci
# cd") and removes the now-empty trailing paragraph.

# Step 1: remove the bookmark from its current location ("Sous-total ").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Step 2: replace the last (empty) paragraph with the new paragraph,
# including the spell-check markers around "Tets" and the relocated
# bookmark, using InsertXML so no extra trailing paragraph is created.
$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastParagraph.Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>Tets</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> pour la ci cd</w:t></w:r>' +
       '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
       '<w:bookmarkEnd w:id="0"/>' +
       '</w:p>'
$r.InsertXML($xml)
